$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A must remain plain text (not auto-converted to a date serial number),
# and must not pick up an extra cell style like the header row uses.
$ws.Range("A94").NumberFormat = "@"
$ws.Range("A94").Value = "02/26/2026"
$ws.Range("A94").Style = "Normal"

$ws.Range("B94").Value = 9766.27
$ws.Range("C94").Value = 0.2387191444856278
$ws.Range("D94").Value = 0.7612808555143722
$ws.Range("E94").Value = -295.48
$ws.Range("F94").Value = -31.6
$ws.Range("G94").Value = -23556.88
$ws.Range("H94").Value = -76.01000000000001
$ws.Range("I94").Value = -1121.45
$ws.Range("J94").Value = -32.48
$ws.Range("K94").Value = -24678.33
$ws.Range("L94").Value = -71.65000000000001
